$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> D-Wave Quantum Inc. / QBTS
$ws.Range("B2").Value = "D-Wave Quantum Inc."
$ws.Range("C2").Value = "QBTS"
$ws.Range("D2").Value = 22.5
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = -2.64
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 80
$ws.Range("I2").Value = 86
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 57.2
$ws.Range("N2").Value = 66.04328690552585

# Row 3 -> Rigetti Computing, Inc. / RGTI
$ws.Range("B3").Value = "Rigetti Computing, Inc."
$ws.Range("C3").Value = "RGTI"
$ws.Range("D3").Value = 23.88
$ws.Range("E3").Value = 29.7
$ws.Range("F3").Value = -10.12
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 70
$ws.Range("I3").Value = 83
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 56
$ws.Range("N3").Value = 66.04328690552585

# Row 4 -> International Business Machines / IBM
$ws.Range("B4").Value = "International Business Machines"
$ws.Range("C4").Value = "IBM"
$ws.Range("D4").Value = 301.78
$ws.Range("E4").Value = 39.8
$ws.Range("F4").Value = -0.77
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 73
$ws.Range("I4").Value = 66
$ws.Range("J4").Value = 63
$ws.Range("K4").Value = 52.2
$ws.Range("N4").Value = 66.04328690552585

# Row 5 -> IonQ, Inc. / IONQ (unchanged name/ticker, values refreshed)
$ws.Range("B5").Value = "IonQ, Inc."
$ws.Range("C5").Value = "IONQ"
$ws.Range("D5").Value = 46.93
$ws.Range("E5").Value = 38.2
$ws.Range("F5").Value = 0.36
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 46
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 63
$ws.Range("K5").Value = 48.8
$ws.Range("N5").Value = 66.04328690552585
